$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "NBS" / source-citation rows (21-22); their content
#     moves further down the sheet (to rows 27-28) to make room for the
#     new "Number of employees / Assets / Turnover" breakdown table.
$ws.Range("A21").Clear()
$ws.Range("A22").Clear()

# --- New sub-header row (bold "title" style, like the row 9 header)
$ws.Range("B16").Value = "Number of employees"
$ws.Range("B16").Font.Bold = $true
$ws.Range("C16").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C16").Font.Bold = $true
$ws.Range("D16").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D16").Font.Bold = $true

# --- New MSME size-class breakdown rows (plain/Normal style)
$ws.Range("A17").Value = "Micro"
$ws.Range("B17").Value = "1-4"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""

$ws.Range("A18").Value = "Small"
$ws.Range("B18").Value = "5-9"
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""

$ws.Range("A19").Value = "Medium"
$ws.Range("B19").Value = "10-49"
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""

$ws.Range("A20").Value = "Large"
$ws.Range("B20").Value = ">49"
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""

# --- Source citation, now living at rows 27-28
$ws.Range("A27").Value = "NBS"
$ws.Range("A27").Font.Bold = $true

$ws.Range("A28").Value = "National Bureau of Statistics (NBS), ""Statistical Yearbook for Southern Sudan"", 2010, p. 142. Available at http://ssnbs.org/statistical-year-book/"
$ws.Range("A28").Font.Italic = $true
